# Scheduled-runner update: refresh computed market-price / profit columns
# (H..N) across several rows on each of the eight item-category sheets.
# Values below are the recomputed figures from the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1688
$ws.Range("J41").Value = 632.25
$ws.Range("L41").Value = 632.25
$ws.Range("N41").Value = -1512.25

$ws.Range("H76").Value = 5847.75
$ws.Range("I76").Value = 5531.6665
$ws.Range("J76").Value = 6796
$ws.Range("K76").Value = 5531.6665
$ws.Range("L76").Value = 6796
$ws.Range("M76").Value = -5216.6665
$ws.Range("N76").Value = -7426

$ws.Range("H79").Value = 5847.75
$ws.Range("I79").Value = 5531.6665
$ws.Range("J79").Value = 6796
$ws.Range("K79").Value = 5531.6665
$ws.Range("L79").Value = 6796
$ws.Range("M79").Value = -4439.6665
$ws.Range("N79").Value = -8980

$ws.Range("H98").Value = 7235
$ws.Range("I98").Value = 7235
$ws.Range("K98").Value = 7235
$ws.Range("M98").Value = -5737

$ws.Range("H106").Value = 75734.92999999999
$ws.Range("I106").Value = 95826.27
$ws.Range("K106").Value = 95826.27
$ws.Range("M106").Value = -95195.27

$ws.Range("H116").Value = 7641.353
$ws.Range("I116").Value = 6941.2104
$ws.Range("K116").Value = 6941.2104
$ws.Range("M116").Value = -3499.2104

$ws.Range("H122").Value = 7235
$ws.Range("I122").Value = 7235
$ws.Range("K122").Value = 21705
$ws.Range("M122").Value = -19255

$ws.Range("H132").Value = 2293.4092
$ws.Range("I132").Value = 2202.8
$ws.Range("K132").Value = 6608.400000000001
$ws.Range("M132").Value = -4078.400000000001

$ws.Range("H137").Value = 2137.25
$ws.Range("I137").Value = 2156.8572
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 6470.571599999999
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -3920.571599999999
$ws.Range("N137").Value = -11100

$ws.Range("H138").Value = 2329.1396
$ws.Range("J138").Value = 3006.2144
$ws.Range("L138").Value = 9018.643199999999
$ws.Range("N138").Value = -19298.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3703.5532
$ws.Range("I32").Value = 4070.162
$ws.Range("J32").Value = 2347.1
$ws.Range("K32").Value = 4070.162
$ws.Range("L32").Value = 2347.1
$ws.Range("M32").Value = -3783.162
$ws.Range("N32").Value = -2921.1

$ws.Range("H122").Value = 1521.1923
$ws.Range("I122").Value = 1437.55
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 4312.65
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -1862.65
$ws.Range("N122").Value = -10300

$ws.Range("H132").Value = 3542
$ws.Range("I132").Value = 3871.5908
$ws.Range("K132").Value = 11614.7724
$ws.Range("M132").Value = -9084.7724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 116749.75
$ws.Range("J76").Value = 116749.75
$ws.Range("L76").Value = 116749.75
$ws.Range("N76").Value = -117379.75

$ws.Range("H79").Value = 116749.75
$ws.Range("J79").Value = 116749.75
$ws.Range("L79").Value = 116749.75
$ws.Range("N79").Value = -118933.75

$ws.Range("H134").Value = 7303.125
$ws.Range("I134").Value = 4090.7
$ws.Range("K134").Value = 12272.1
$ws.Range("M134").Value = -9737.099999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 201766.2
$ws.Range("I62").Value = 251710.75
$ws.Range("K62").Value = 251710.75
$ws.Range("M62").Value = -251086.75

$ws.Range("H65").Value = 201766.2
$ws.Range("I65").Value = 251710.75
$ws.Range("K65").Value = 1258553.75
$ws.Range("M65").Value = -1255433.75

$ws.Range("H107").Value = 35716370
$ws.Range("J107").Value = 3348.2
$ws.Range("L107").Value = 3348.2
$ws.Range("N107").Value = -7188.2

$ws.Range("H132").Value = 2356.1765
$ws.Range("I132").Value = 2293.0312
$ws.Range("J132").Value = 3366.5
$ws.Range("K132").Value = 6879.0936
$ws.Range("L132").Value = 10099.5
$ws.Range("M132").Value = -4349.0936
$ws.Range("N132").Value = -15159.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2095.3333
$ws.Range("I14").Value = 2095.3333
$ws.Range("K14").Value = 6285.999899999999
$ws.Range("M14").Value = -6112.999899999999

$ws.Range("H23").Value = 682
$ws.Range("I23").Value = 667.125
$ws.Range("J23").Value = 701.8333
$ws.Range("K23").Value = 2001.375
$ws.Range("L23").Value = 2105.4999
$ws.Range("M23").Value = -1766.375
$ws.Range("N23").Value = -2575.4999

$ws.Range("H68").Value = 1369.0625
$ws.Range("J68").Value = 1616.2941
$ws.Range("L68").Value = 4848.8823
$ws.Range("N68").Value = -6470.8823

$ws.Range("H71").Value = 1369.0625
$ws.Range("J71").Value = 1616.2941
$ws.Range("L71").Value = 14546.6469
$ws.Range("N71").Value = -22658.6469

$ws.Range("H117").Value = 302498
$ws.Range("I117").Value = 658.7778
$ws.Range("J117").Value = 496537.5
$ws.Range("K117").Value = 1976.3334
$ws.Range("L117").Value = 1489612.5
$ws.Range("M117").Value = 1465.6666
$ws.Range("N117").Value = -1496496.5

$ws.Range("H129").Value = 3932.8
$ws.Range("I129").Value = 1287.5
$ws.Range("J129").Value = 5696.3335
$ws.Range("K129").Value = 3862.5
$ws.Range("L129").Value = 17089.0005
$ws.Range("M129").Value = 1137.5
$ws.Range("N129").Value = -27089.0005

$ws.Range("H134").Value = 5615.7393
$ws.Range("I134").Value = 5219.0527
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 15657.1581
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = -10587.1581
$ws.Range("N134").Value = -32640

$ws.Range("H136").Value = 4117.5713
$ws.Range("I136").Value = 4117.5713
$ws.Range("K136").Value = 12352.7139
$ws.Range("M136").Value = -7252.713899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1582210.6
$ws.Range("I132").Value = 1669889
$ws.Range("K132").Value = 5009667
$ws.Range("M132").Value = -5007137

$ws.Range("H134").Value = 38332.43
$ws.Range("J134").Value = 38332.43
$ws.Range("L134").Value = 114997.29
$ws.Range("N134").Value = -120067.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3992.2727
$ws.Range("I7").Value = 3507.75
$ws.Range("K7").Value = 3507.75
$ws.Range("M7").Value = -3395.75

$ws.Range("H32").Value = 14658.223
$ws.Range("I32").Value = 7417.7144
$ws.Range("K32").Value = 7417.7144
$ws.Range("M32").Value = -7100.7144

$ws.Range("H40").Value = 5902.9375
$ws.Range("I40").Value = 5561.4614
$ws.Range("K40").Value = 5561.4614
$ws.Range("M40").Value = -5425.4614

$ws.Range("H122").Value = 3529.111
$ws.Range("I122").Value = 3408.1667
$ws.Range("K122").Value = 10224.5001
$ws.Range("M122").Value = -7774.500100000001

$ws.Range("H126").Value = 3992.2727
$ws.Range("I126").Value = 3507.75
$ws.Range("K126").Value = 10523.25
$ws.Range("M126").Value = -8053.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 45024

$ws.Range("H43").Value = 41015
$ws.Range("J43").Value = 41015
$ws.Range("L43").Value = 41015
$ws.Range("N43").Value = -41313

$ws.Range("H50").Value = 40042
$ws.Range("J50").Value = 40042
$ws.Range("L50").Value = 40042
$ws.Range("N50").Value = -41304

$ws.Range("H53").Value = 29650.4
$ws.Range("J53").Value = 29650.4
$ws.Range("L53").Value = 29650.4
$ws.Range("N53").Value = -30864.4

$ws.Range("H96").Value = 1995.75
$ws.Range("J96").Value = 2015.75
$ws.Range("L96").Value = 2015.75
$ws.Range("N96").Value = -4761.75

$ws.Range("H100").Value = 40000596
$ws.Range("I100").Value = 62500612
$ws.Range("J100").Value = 572
$ws.Range("K100").Value = 125001224
$ws.Range("L100").Value = 1144
$ws.Range("M100").Value = -125000683
$ws.Range("N100").Value = -2226

$ws.Range("H132").Value = 3477.7812
$ws.Range("I132").Value = 3510
$ws.Range("J132").Value = 3166.3333
$ws.Range("K132").Value = 10530
$ws.Range("L132").Value = 9498.999899999999
$ws.Range("M132").Value = -8000
$ws.Range("N132").Value = -14558.9999

$ws.Range("H136").Value = 3346214.8
$ws.Range("J136").Value = 1988.2222
$ws.Range("L136").Value = 5964.6666
$ws.Range("N136").Value = -11064.6666
